$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 154 ---
$ws.Range("A153").Copy()
$ws.Range("A154").PasteSpecial(-4122)
$ws.Range("A154").Value = 45504.2916666667

$ws.Range("B154").Value = 0
$ws.Range("C154").Value = 2.94000005722046
$ws.Range("D154").Value = 2.94000005722046
$ws.Range("E154").Value = 2.94000005722046
$ws.Range("F154").Value = 2.94000005722046

$ws.Range("G154").NumberFormat = "@"
$ws.Range("G154").Value = "2.94000005722046"
$ws.Range("G154").ClearFormats()

$ws.Range("H154").NumberFormat = "@"
$ws.Range("H154").Value = "AGAIN.MI"
$ws.Range("H154").ClearFormats()

# --- Row 155 ---
$ws.Range("A153").Copy()
$ws.Range("A155").PasteSpecial(-4122)
$ws.Range("A155").Value = 45505.6319444444

$ws.Range("B155").Value = 5000
$ws.Range("C155").Value = 2.98000001907349
$ws.Range("D155").Value = 2.90000009536743
$ws.Range("E155").Value = 2.90000009536743
$ws.Range("F155").Value = 2.94000005722046

$ws.Range("G155").NumberFormat = "@"
$ws.Range("G155").Value = "2.94000005722046"
$ws.Range("G155").ClearFormats()

$ws.Range("H155").NumberFormat = "@"
$ws.Range("H155").Value = "AGAIN.MI"
$ws.Range("H155").ClearFormats()

$excel.CutCopyMode = $false
